# "add constructing of Student object"
#
# The "Grades comments" worksheet is removed; a per-grade "Comment" column
# is added directly onto the "Grades" worksheet instead (column E).

$wb = $excel.ActiveWorkbook

$gradesWs = $wb.Worksheets("Grades")

# Drop the now-redundant sheet.
$wb.Worksheets("Grades comments").Delete()

# Add the "Comment" column onto "Grades", one remark per grade row.
$gradesWs.Range("E1").Value = "Comment"
$gradesWs.Range("E2").Value = "No comment."
$gradesWs.Range("E3").Value = "Bad job :("
$gradesWs.Range("E4").Value = "No comment."
$gradesWs.Range("E5").Value = "Aweful work"
$gradesWs.Range("E6").Value = "Well done!"

$gradesWs.Columns("E").ColumnWidth = 11.92

# Move the active tab/selection onto "Grades".
$gradesWs.Activate()
$gradesWs.Range("E7").Select()
